# Powerpoint writer: consolidate text run nodes.
# Merge the leading "First"/" " runs (and "Third"/" " runs) into a single
# run per title, since both runs share identical (empty) run properties.
# Re-assigning the text of a sub-range spanning exactly those runs causes
# the writer to collapse them into one <a:r> while leaving the trailing
# "slide" run untouched.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "First "

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 6).Text = "Third "
